$d = $word.ActiveDocument

# Locate the first "TODO" placeholder (the one right after the
# "GitHub repo with all workshop files" line) and turn it into a
# hyperlink pointing at the workshop's GitHub repo.
$url = "https://github.com/RGreinacher/geospatial-big-data-fasrc"

$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("TODO", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $d.Hyperlinks.Add($rng, $url, $null, $null, $url) | Out-Null
}
